# Horarios Línea 141 update (816) - applies the scrape refresh at 04:01:13
$wb = $excel.ActiveWorkbook

$newTime = "04:01:13"

# ----------------------------------------------------------------------
# Sheet 1: LP1912
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 19"

# Row 14: only the scrape timestamp and the minutes-elapsed column changed
$ws1.Cells.Item(14, 1).Value = $newTime
$ws1.Cells.Item(14, 4).Value = 0

# Rows 16-19: refreshed scrape data (each arrival shifted up one stop)
$ws1.Cells.Item(16, 1).Value = $newTime
$ws1.Cells.Item(16, 2).Value = "04:46"
$ws1.Cells.Item(16, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(16, 4).Value = 45
$ws1.Cells.Item(16, 5).Value = "LP1912"

$ws1.Cells.Item(17, 1).Value = $newTime
$ws1.Cells.Item(17, 2).Value = "04:53"
$ws1.Cells.Item(17, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(17, 4).Value = 52
$ws1.Cells.Item(17, 5).Value = "LP1912"

$ws1.Cells.Item(18, 1).Value = $newTime
$ws1.Cells.Item(18, 2).Value = "05:16"
$ws1.Cells.Item(18, 3).Value = "17_ROMERO"
$ws1.Cells.Item(18, 4).Value = 75
$ws1.Cells.Item(18, 5).Value = "LP1912"

$ws1.Cells.Item(19, 1).Value = $newTime
$ws1.Cells.Item(19, 2).Value = "05:22"
$ws1.Cells.Item(19, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(19, 4).Value = 81
$ws1.Cells.Item(19, 5).Value = "LP1912"

# New rows 20-24 appended with the latest scrape
$ws1.Cells.Item(20, 1).Value = "03:35:49"
$ws1.Cells.Item(20, 2).Value = "05:28"
$ws1.Cells.Item(20, 3).Value = "14_ABASTO"
$ws1.Cells.Item(20, 4).Value = 113
$ws1.Cells.Item(20, 5).Value = "LP1912"

$ws1.Cells.Item(21, 1).Value = "03:35:49"
$ws1.Cells.Item(21, 2).Value = "05:34"
$ws1.Cells.Item(21, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(21, 4).Value = 119
$ws1.Cells.Item(21, 5).Value = "LP1912"

$ws1.Cells.Item(22, 1).Value = $newTime
$ws1.Cells.Item(22, 2).Value = "05:35"
$ws1.Cells.Item(22, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(22, 4).Value = 94
$ws1.Cells.Item(22, 5).Value = "LP1912"

$ws1.Cells.Item(23, 1).Value = $newTime
$ws1.Cells.Item(23, 2).Value = "05:37"
$ws1.Cells.Item(23, 3).Value = "14_ABASTO"
$ws1.Cells.Item(23, 4).Value = 96
$ws1.Cells.Item(23, 5).Value = "LP1912"

$ws1.Cells.Item(24, 1).Value = $newTime
$ws1.Cells.Item(24, 2).Value = "05:46"
$ws1.Cells.Item(24, 3).Value = "15_ABASTO"
$ws1.Cells.Item(24, 4).Value = 105
$ws1.Cells.Item(24, 5).Value = "LP1912"

# ----------------------------------------------------------------------
# Sheet 2: LP1912-215
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 8"

# A new arrival (215A_EL PATO) is inserted before the old row 11, pushing
# the existing 215B_EL PATO row down to row 12.
$ws2.Rows.Item(11).Insert()

$ws2.Cells.Item(11, 1).Value = $newTime
$ws2.Cells.Item(11, 2).Value = "04:46"
$ws2.Cells.Item(11, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(11, 4).Value = 45
$ws2.Cells.Item(11, 5).Value = "LP1912"

# New row 13 appended with the latest scrape
$ws2.Cells.Item(13, 1).Value = $newTime
$ws2.Cells.Item(13, 2).Value = "05:35"
$ws2.Cells.Item(13, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(13, 4).Value = 94
$ws2.Cells.Item(13, 5).Value = "LP1912"

# ----------------------------------------------------------------------
# Sheet 3: 6203-6173
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
$ws3.Range("A3").Value = "Total filas: 2"

# New row 7 appended with the latest scrape
$ws3.Cells.Item(7, 1).Value = $newTime
$ws3.Cells.Item(7, 2).Value = "05:44"
$ws3.Cells.Item(7, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(7, 4).Value = 103
$ws3.Cells.Item(7, 5).Value = "L6173"
